$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 ("I0") and J1 ("IF") should look like the existing
# header cells (e.g. H1): copy H1's formatting onto I1:J1, then set values.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for the new I0 / IF columns (rows 2-11)
$values = @(
    @(8, 8),
    @(6, 6),
    @(8, 8),
    @(6, 6),
    @(10, 10),
    @(6, 6),
    @(6, 7),
    @(5, 5),
    @(4, 5),
    @(7, 7)
)

for ($i = 0; $i -lt $values.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $values[$i][0]
    $ws.Cells.Item($row, 10).Value = $values[$i][1]
}
